$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misread as a plain number by Excel
# (single-dot decimal-looking strings). Force Text format first so the
# literal string (incl. any trailing zero) survives, then restore the
# default "Normal" style so no stray formatting is left behind.
$textForceRows = @(5, 6, 7, 9, 10, 13, 14, 16, 20, 22, 23, 24, 27, 28, 30, 31, 35, 39, 40, 43, 44, 47, 48, 50, 51)
foreach ($r in $textForceRows) {
  $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(5, 4).Value = '253.23'
$ws.Cells.Item(6, 4).Value = '0.653'
$ws.Cells.Item(7, 4).Value = '65.81'
$ws.Cells.Item(9, 4).Value = '0.384'
$ws.Cells.Item(10, 4).Value = '59.56'
$ws.Cells.Item(13, 4).Value = '0.925'
$ws.Cells.Item(14, 4).Value = '14.88'
$ws.Cells.Item(16, 4).Value = '22.37'
$ws.Cells.Item(20, 4).Value = '73.64'
$ws.Cells.Item(22, 4).Value = '5.51'
$ws.Cells.Item(23, 4).Value = '239.79'
$ws.Cells.Item(24, 4).Value = '2.70'
$ws.Cells.Item(27, 4).Value = '9.98'
$ws.Cells.Item(28, 4).Value = '162.07'
$ws.Cells.Item(30, 4).Value = '0.124'
$ws.Cells.Item(31, 4).Value = '5.26'
$ws.Cells.Item(35, 4).Value = '0.0627'
$ws.Cells.Item(39, 4).Value = '6.04'
$ws.Cells.Item(40, 4).Value = '3.04'
$ws.Cells.Item(43, 4).Value = '1.24'
$ws.Cells.Item(44, 4).Value = '17.78'
$ws.Cells.Item(47, 4).Value = '97.20'
$ws.Cells.Item(48, 4).Value = '7.92'
$ws.Cells.Item(50, 4).Value = '3.90'
$ws.Cells.Item(51, 4).Value = '2.95'

foreach ($r in $textForceRows) {
  $ws.Cells.Item($r, 4).Style = "Normal"
}

# Remaining cells: plain text (URLs, coin names) or percentage strings that
# are never ambiguous with a numeric literal (two embedded dots, spaces, % sign).
$ws.Cells.Item(2, 4).Value = '37.461.12'
$ws.Cells.Item(2, 5).Value = '  +5.14%  '
$ws.Cells.Item(3, 4).Value = '2.054.81'
$ws.Cells.Item(3, 5).Value = '  +3.53%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 5).Value = '  +3.12%  '
$ws.Cells.Item(6, 5).Value = '  +2.21%  '
$ws.Cells.Item(7, 5).Value = '  +12.96%  '
$ws.Cells.Item(8, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 5).Value = '  +6.48%  '
$ws.Cells.Item(10, 5).Value = '  +0.95%  '
$ws.Cells.Item(11, 5).Value = '  +4.70%  '
$ws.Cells.Item(12, 5).Value = '  +1.38%  '
$ws.Cells.Item(13, 5).Value = '  -2.66%  '
$ws.Cells.Item(14, 5).Value = '  +2.36%  '
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value = '2.356.70'
$ws.Cells.Item(15, 5).Value = '  +3.67%  '
$ws.Cells.Item(16, 2).Value = 'Avalanche'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(16, 5).Value = '  +25.10%  '
$ws.Cells.Item(17, 5).Value = '  +5.03%  '
$ws.Cells.Item(18, 4).Value = '2.053.08'
$ws.Cells.Item(18, 5).Value = '  +3.52%  '
$ws.Cells.Item(19, 4).Value = '37.345.84'
$ws.Cells.Item(19, 5).Value = '  +4.88%  '
$ws.Cells.Item(20, 5).Value = '  +3.13%  '
$ws.Cells.Item(21, 5).Value = '  +3.62%  '
$ws.Cells.Item(22, 5).Value = '  +6.39%  '
$ws.Cells.Item(23, 5).Value = '  +2.92%  '
$ws.Cells.Item(24, 5).Value = '  +4.65%  '
$ws.Cells.Item(25, 5).Value = '  -0.09%  '
$ws.Cells.Item(26, 5).Value = '  +5.35%  '
$ws.Cells.Item(27, 5).Value = '  +8.83%  '
$ws.Cells.Item(28, 5).Value = '  -1.66%  '
$ws.Cells.Item(29, 5).Value = '  +4.04%  '
$ws.Cells.Item(30, 5).Value = '  +28.48%  '
$ws.Cells.Item(31, 5).Value = '  +8.14%  '
$ws.Cells.Item(32, 5).Value = '  +2.37%  '
$ws.Cells.Item(33, 5).Value = '  +9.65%  '
$ws.Cells.Item(34, 5).Value = '  +8.25%  '
$ws.Cells.Item(35, 5).Value = '  +5.57%  '
$ws.Cells.Item(36, 5).Value = '  +1.18%  '
$ws.Cells.Item(37, 5).Value = '  +0.06%  '
$ws.Cells.Item(38, 5).Value = '  +4.01%  '
$ws.Cells.Item(39, 5).Value = '  +14.94%  '
$ws.Cells.Item(40, 5).Value = '  +35.22%  '
$ws.Cells.Item(41, 5).Value = '  +16.26%  '
$ws.Cells.Item(42, 5).Value = '  +4.37%  '
$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43, 5).Value = '  +1.43%  '
$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(44, 5).Value = '  +9.54%  '
$ws.Cells.Item(45, 5).Value = '  +5.91%  '
$ws.Cells.Item(47, 5).Value = '  +5.24%  '
$ws.Cells.Item(48, 5).Value = '  +2.50%  '
$ws.Cells.Item(49, 4).Value = '1.417.96'
$ws.Cells.Item(49, 5).Value = '  +3.18%  '
$ws.Cells.Item(50, 2).Value = 'FTXToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(50, 5).Value = '  +12.51%  '
$ws.Cells.Item(51, 2).Value = 'MXToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(51, 5).Value = '  +1.60%  '

Write-Output "Applied 92 cell updates (25 text-forced, 67 plain)"
